$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new observation ("2026/02/26", "木", 9, 16) was recorded and needs to be
# inserted in date order at row 861, pushing the existing rows 861:902 down
# to 862:903 (dimension grows from D902 to D903).
$ws.Rows("861").Insert()

# Force column A to stay plain text ("YYYY/MM/DD" string), matching every
# other row in the date column, instead of being auto-parsed into a date
# serial number.
$ws.Range("A861").NumberFormat = "@"
$ws.Range("A861").Value = "2026/02/26"
$ws.Range("B861").Value = "木"
$ws.Range("C861").Value = 9
$ws.Range("D861").Value = 16
